# Renames the inline "Picture" objects that Word shows for the three
# header/footer logo images. The two Pearson-logo pictures (one in each
# footer) are renamed from "image2.png" to "image1.png", and the BTEC
# logo picture (in the first-page header) is renamed from "image1.jpg"
# to "image2.jpg".
#
# InlineShape.Name cannot be set reliably straight off a
# Header/Footer.Range.InlineShapes collection item in this host, so each
# shape's own Range is selected first and then re-fetched through
# $word.Selection.InlineShapes - that re-seated handle accepts the
# rename.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($hostRange, $newName) {
    if ($hostRange.Exists -and $hostRange.Range.InlineShapes.Count -gt 0) {
        $shp = $hostRange.Range.InlineShapes.Item(1)
        $shp.Range.Select()
        $selShp = $word.Selection.InlineShapes.Item(1)
        $selShp.Name = $newName
    }
}

# Pearson Edexcel logo - default (primary) footer
Rename-InlineLogo $sec.Footers.Item(1) "image1.png"

# Pearson Edexcel logo - first-page footer
Rename-InlineLogo $sec.Footers.Item(2) "image1.png"

# BTEC logo - first-page header
Rename-InlineLogo $sec.Headers.Item(2) "image2.jpg"
